# Update "想去人数" (number of people interested) figures to reflect
# newly generated output data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 980
$wsExhibit.Range("F3").Value = 1982
$wsExhibit.Range("F4").Value = 442

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 980
$wsAll.Range("F5").Value = 1982
$wsAll.Range("F6").Value = 442
